$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New raw data rows (101-108) -------------------------------------------------
# Columns: A Sampled(date) B ID C Volume(mL) D Tray weight(g) E Desicator(date)
#          F Desicator(g) G Furance(date) H Furnace(g) I Analyzer J POC(mg) K mg/L

$data = @(
    @(45526, 5,    1060, 1.0873999999999999, 45534, 1.3194999999999999, 45538, 1.3109999999999999, "SH"),
    @(45527, 9,    556,  1.1223000000000001, 45534, 1.339,               45538, 1.3368,              "SH"),
    @(45527, 15,   558,  1.1051,             45534, 1.3209,              45538, 1.319,               "SH"),
    @(45527, 7,    559,  1.0774999999999999, 45534, 1.2995000000000001, 45538, 1.2928999999999999,  "SH"),
    @(45526, "5a", 1077, 1.0660000000000001, 45534, 1.3069,              45538, 1.2948,              "SH"),
    @(45527, "6a", 572,  1.1362000000000001, 45534, 1.3665,              45538, 1.3612,              "SH"),
    @(45527, 13,   559,  1.1144000000000001, 45534, 1.3347,              45538, 1.3326,              "SH"),
    @(45527, 3,    559,  1.1178999999999999, 45534, 1.3354999999999999, 45538, 1.3354999999999999,  "SH")
)

$startRow = 101
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]   # A - Sampled date
    $ws.Cells.Item($r, 2).Value = $row[1]   # B - ID
    $ws.Cells.Item($r, 3).Value = $row[2]   # C - Volume (mL)
    $ws.Cells.Item($r, 4).Value = $row[3]   # D - Tray weight (g)
    $ws.Cells.Item($r, 5).Value = $row[4]   # E - Desicator date
    $ws.Cells.Item($r, 6).Value = $row[5]   # F - Desicator (g)
    $ws.Cells.Item($r, 7).Value = $row[6]   # G - Furance date
    $ws.Cells.Item($r, 8).Value = $row[7]   # H - Furnace (g)
    $ws.Cells.Item($r, 9).Value = $row[8]   # I - Analyzer
}

# --- Copy number formats from the prior row (100) so date cells keep the same styles
$ws.Range("A100").Copy()
$ws.Range("A101:A108").PasteSpecial(-4122)

$ws.Range("E100").Copy()
$ws.Range("E101:E108").PasteSpecial(-4122)

$ws.Range("A100").Copy()
$ws.Range("G101:G108").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Extend the shared formulas down through row 108 ----------------------------
$ws.Range("J101:J108").Formula = "=(F101-H101)*1000"
$ws.Range("K101:K108").Formula = "=J101/(C101/1000)"

# --- Update sheet view (selection) -----------------------------------------------
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H98").Select()
